$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell in H1: same text + same style (s="1") as the rest of row 1.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# New "Save" column values for rows 2-9 (plain numbers, no special style).
$saveValues = @(1, 1, 0, 0, 0, 0, 0, 0)
for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
